# This script updates the "Recorded By" column (column G) on the
# "Session Analysis Results" sheet, reordering the comma-separated list of
# recorder names/emails for specific known values, matching the upstream
# sync from the main repository.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact-value mapping: old "Recorded By" text -> new "Recorded By" text.
# Only cells whose current value exactly matches one of these keys are
# touched; everything else (single-name cells, or combinations not part
# of this sync) is left untouched.
$map = @{
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

# Determine the used range so we know how many rows to scan.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $current = $cell.Value()

    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
